$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 grew slightly taller (18.75 -> 19.5), matching the header row.
$ws.Rows.Item(2).RowHeight = 19.5

# C2/D2 picked up the same explicit black font color already used by the
# header cells in row 1 (C1/D1), instead of the automatic/theme color.
$ws.Range("C2:D2").Font.Color = 0

# Balance went up.
$ws.Range("D2").Value = 1891762

# Bearer token was refreshed with a new tokenId/iat.
$ws.Range("E2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VySWQiOiI0ODI5OGVhMC0yNDBhLTExZWUtOWMwNC1iMzcyMDk2MTViOGIiLCJhcHBJZCI6MjMyLCJ0b2tlbklkIjoiNWVmMWNjMWQtZGNjYi00OGUxLThmODItOGUxOTMxZGVkODU3IiwiaWF0IjoxNzE1MDUwNzg4fQ.1ND80rCzorESlLsqr3bOdAuquHBAbMY_nV1Yn6hFjpU"
